$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("C3").Value = "TERM BLOCK HDR 2POS 90DEG 3.5MM"
$ws.Range("D3").Value = "277-2416-ND"
$ws.Range("E3").Value = "Phoenix Contact"
$ws.Range("F3").Value = 1844210
$ws.Range("H3").Value = 1.42

# Row 4
$ws.Range("C4").Value = "CONN HEADER R/A 10POS 2.54MM"
$ws.Range("D4").Value = "2057-PH1RB-10-UA-ND"
$ws.Range("E4").Value = "Adam Tech"
$ws.Range("F4").Value = "PH1RB-10-UA"
$ws.Range("H4").Value = 0.41

# Row 5
$ws.Range("C5").Value = "CONN HEADER VERT 24POS 2.54MM"
$ws.Range("D5").Value = "S1131EC-24-ND"
$ws.Range("E5").Value = "Sullins Connector Solutions"
$ws.Range("F5").Value = "PRPC024SACN-RC"
$ws.Range("H5").Value = 1.1

# Row 6
$ws.Range("C6").Value = "CONN RCPT USB2.0 MICRO B SMD R/A"
$ws.Range("D6").Value = "609-4613-1-ND"
$ws.Range("E6").Value = "Amphenol ICC (FCI)"
$ws.Range("F6").Value = "10118192-0001LF"
$ws.Range("H6").Value = 0.66

# Row 7
$ws.Range("C7").Value = "IC REG LINEAR 5V 1A SOT223"
$ws.Range("D7").Value = "NCP1117ST50T3GOSCT-ND"
$ws.Range("E7").Value = "ON Semiconductor"
$ws.Range("F7").Value = "NCP1117ST50T3G"
$ws.Range("H7").Value = 0.76

# Row 8
$ws.Range("C8").Value = "IC REG LINEAR 3.3V 1.2A SOT223"
$ws.Range("D8").Value = "497-17239-1-ND"
$ws.Range("E8").Value = "STMicroelectronics"
$ws.Range("F8").Value = "LDL1117S33R"
$ws.Range("H8").Value = 0.69

# Row 9
$ws.Range("C9").Value = "IC MCU 32BIT 1MB FLASH 64TQFP"
$ws.Range("D9").Value = "ATSAMD51J20A-AU-ND"
$ws.Range("E9").Value = "Microchip Technology"
$ws.Range("F9").Value = "ATSAMD51J20A-AU"
$ws.Range("H9").Value = 7.89

# Row 10
$ws.Range("C10").Value = "IC FLASH 4M SPI 104MHZ 8SOIC"
$ws.Range("D10").Value = "W25X40CLSNIG-ND"
$ws.Range("E10").Value = "Winbond Electronics"
$ws.Range("F10").Value = "W25X40CLSNIG"
$ws.Range("H10").Value = 0.59

# Row 11
$ws.Range("C11").Value = "DIODE SCHOTTKY 40V 1A SOD123"
$ws.Range("D11").Value = "1N5819HW-FDICT-ND"
$ws.Range("E11").Value = "Diodes Incorporated"
$ws.Range("F11").Value = "1N5819HW-7-F"
$ws.Range("H11").Value = 0.61

# Row 12
$ws.Range("C12").Value = "LED GREEN DIFFUSED 0603 SMD"
$ws.Range("D12").Value = "475-3118-1-ND"
$ws.Range("E12").Value = "OSRAM Opto Semiconductors Inc."
$ws.Range("F12").Value = "LG L29K-F2J1-24-Z"
$ws.Range("H12").Value = 0.62

# Row 13
$ws.Range("C13").Value = "CRYSTAL 32.768KHZ 12.5PF SMD"
$ws.Range("D13").Value = "2195-CM7V-T1A-32.768KHZ-12.5PF-20PPM-TA-QCCT-ND"
$ws.Range("E13").Value = "Micro Crystal AG"
$ws.Range("F13").Value = "CM7V-T1A-32.768KHZ-12.5PF-20PPM-TA-QC"
$ws.Range("H13").Value = 0.67

# Row 14
$ws.Range("C14").Value = "TACT 4.5 X 4.5, 3.8 MM H, 2.5N,"
$ws.Range("D14").Value = "PTS647SK38SMTR2LFSCT-ND"
$ws.Range("E14").Value = "C&K"
$ws.Range("F14").Value = "PTS 647 SK38 SMTR2 LFS"
$ws.Range("H14").Value = 0.23

# Row 15
$ws.Range("C15").Value = "FIXED IND 10UH 150MA 360 MOHM"
$ws.Range("D15").Value = "587-2045-1-ND"
$ws.Range("E15").Value = "Taiyo Yuden"
$ws.Range("F15").Value = "LBR2012T100K"
$ws.Range("H15").Value = 0.21

# Row 16
$ws.Range("C16").Value = "CAP CER 22UF 25V X5R 0805"
$ws.Range("D16").Value = "1276-2908-1-ND"
$ws.Range("E16").Value = "Samsung Electro-Mechanics"
$ws.Range("F16").Value = " CL21A226MAQNNNE"
$ws.Range("H16").Value = 0.8

# Row 17
$ws.Range("C17").Value = "CAP CER 47UF 6.3V X5R 0805"
$ws.Range("D17").Value = "490-9960-1-ND"
$ws.Range("E17").Value = "Murata Electronics"
$ws.Range("F17").Value = "GRM21BR60J476ME15L"
$ws.Range("H17").Value = 0.8

# Row 18
$ws.Range("C18").Value = "CAP CER 15PF 50V C0G/NPO 0402"
$ws.Range("D18").Value = "311-1642-1-ND"
$ws.Range("E18").Value = "Yageo"
$ws.Range("F18").Value = "CC0402FRNPO9BN150"
$ws.Range("H18").Value = 0.18

# Row 19
$ws.Range("C19").Value = "CAP CER 10UF 6.3V X5R 0805"
$ws.Range("D19").Value = "1276-2405-1-ND"
$ws.Range("E19").Value = "Samsung Electro-Mechanics"
$ws.Range("F19").Value = "CL21A106KQCLRNC"
$ws.Range("H19").Value = 0.17

# Row 20
$ws.Range("C20").Value = "CAP CER 0.1UF 25V X7R 0603"
$ws.Range("D20").Value = "311-1341-1-ND"
$ws.Range("E20").Value = "Yageo"
$ws.Range("F20").Value = "CC0603KRX7R8BB104"
$ws.Range("H20").Value = 0.15

# Row 21
$ws.Range("C21").Value = "CAP CER 4.7UF 6.3V X5R 0603"
$ws.Range("D21").Value = "1276-1045-1-ND"
$ws.Range("E21").Value = "Samsung Electro-Mechanics"
$ws.Range("F21").Value = "CL10A475KQ8NNNC"
$ws.Range("H21").Value = 0.15

# Row 22
$ws.Range("C22").Value = "RES 820 OHM 1% 1/8W 0805"
$ws.Range("D22").Value = "RMCF0805FT820RCT-ND"
$ws.Range("E22").Value = "Stackpole Electronics Inc"
$ws.Range("F22").Value = "RMCF0805FT820R"
$ws.Range("H22").Value = 0.16

# Row 23
$ws.Range("C23").Value = "RES SMD 10K OHM 1% 1/8W 0805"
$ws.Range("D23").Value = "311-10.0KCRCT-ND"
$ws.Range("E23").Value = "Yageo"
$ws.Range("F23").Value = "RC0805FR-0710KL"
$ws.Range("H23").Value = 0.16

# Remove now-unused trailing rows (previously blank placeholder rows 36-39)
$ws.Rows("36:39").Delete()

# Update the active selection to match the post-edit state
$ws.Range("C24").Select()
